$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.6753301551942219
$ws.Range("C2").Value = 114.8270160096505
$ws.Range("D2").Value = 0.1575252929769615
$ws.Range("E2").Value = 645.3272768299601
$ws.Range("G2").Value = 760.9871482877818

$ws.Range("B3").Value = 0.003994804209775715
$ws.Range("C3").Value = 1.667794583268128
$ws.Range("D3").Value = 3.900430680208489
$ws.Range("E3").Value = 9353990175.932438
$ws.Range("G3").Value = 9353990181.504658
